$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ExcelModuleDemoToDoItem")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- Sheet2: mirror Sheet1's description/category/subcategory columns (A:C, rows 1-8) ---
for ($r = 1; $r -le 8; $r++) {
    for ($c = 1; $c -le 3; $c++) {
        $srcVal = $ws1.Cells.Item($r, $c).Value()
        if ($srcVal -ne $null) {
            $ws2.Cells.Item($r, $c).Value = $srcVal
        }
    }
    # carry over the styled-but-empty "A" column look used on Sheet1 (font/alignment)
    if ($ws1.Cells.Item($r, 1).Value() -eq $null) {
        $ws1.Range("A" + $r).Copy()
        $ws2.Range("A" + $r).PasteSpecial(-4122)
    }
}

# --- Sheet3: header row + a new "Another Item" row, with the same blank styled rows below ---
$ws3.Cells.Item(1, 1).Value = "category"
$ws3.Cells.Item(1, 2).Value = "subcategory"
$ws3.Cells.Item(1, 3).Value = "description"

$ws3.Cells.Item(4, 1).Value = "Domestic"
$ws3.Cells.Item(4, 2).Value = "Shopping"
$ws3.Cells.Item(4, 3).Value = "Another Item"

$ws1.Range("A9").Copy()
$ws3.Range("A5").PasteSpecial(-4122)
$ws1.Range("A9").Copy()
$ws3.Range("A6").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Sheet3 becomes the active/selected tab (Sheet1 loses tabSelected)
$ws3.Activate()
